# Scheduled market-data refresh: updates currentAveragePrice / LevePrice /
# LeveProfit columns (H:N) for a handful of leve rows across the eight Job
# sheets, as produced by the Coeurl-server price-fetch runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 1201
$ws.Cells.Item(98, 9).Value = 1135.6086
$ws.Cells.Item(98, 11).Value = 1135.6086
$ws.Cells.Item(98, 13).Value = 362.3914

$ws.Cells.Item(122, 8).Value = 1201
$ws.Cells.Item(122, 9).Value = 1135.6086
$ws.Cells.Item(122, 11).Value = 3406.8258
$ws.Cells.Item(122, 13).Value = -956.8258000000001

$ws.Cells.Item(138, 8).Value = 3083.5444
$ws.Cells.Item(138, 9).Value = 1679.0646
$ws.Cells.Item(138, 10).Value = 3821.4915
$ws.Cells.Item(138, 11).Value = 5037.1938
$ws.Cells.Item(138, 12).Value = 11464.4745
$ws.Cells.Item(138, 13).Value = 102.8062
$ws.Cells.Item(138, 14).Value = -21744.4745

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 34485932
$ws.Cells.Item(2, 9).Value = 41669300
$ws.Cells.Item(2, 11).Value = 41669300
$ws.Cells.Item(2, 13).Value = -41669187

$ws.Cells.Item(116, 8).Value = 34485932
$ws.Cells.Item(116, 9).Value = 41669300
$ws.Cells.Item(116, 11).Value = 41669300
$ws.Cells.Item(116, 13).Value = -41667006

$ws.Cells.Item(132, 8).Value = 2970
$ws.Cells.Item(132, 9).Value = 2289.9565
$ws.Cells.Item(132, 10).Value = 4534.1
$ws.Cells.Item(132, 11).Value = 6869.869499999999
$ws.Cells.Item(132, 12).Value = 13602.3
$ws.Cells.Item(132, 13).Value = -4339.869499999999
$ws.Cells.Item(132, 14).Value = -18662.3

$ws.Cells.Item(140, 8).Value = 51138.75
$ws.Cells.Item(140, 9).Value = 40000
$ws.Cells.Item(140, 10).Value = 54851.668
$ws.Cells.Item(140, 11).Value = 40000
$ws.Cells.Item(140, 12).Value = 54851.668
$ws.Cells.Item(140, 13).Value = -34820
$ws.Cells.Item(140, 14).Value = -65211.668

$ws.Cells.Item(141, 8).Value = 45185
$ws.Cells.Item(141, 9).Value = 40000
$ws.Cells.Item(141, 10).Value = 47777.5
$ws.Cells.Item(141, 11).Value = 40000
$ws.Cells.Item(141, 12).Value = 47777.5
$ws.Cells.Item(141, 13).Value = -34820
$ws.Cells.Item(141, 14).Value = -58137.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 34485932
$ws.Cells.Item(3, 9).Value = 41669300
$ws.Cells.Item(3, 11).Value = 41669300
$ws.Cells.Item(3, 13).Value = -41669186

$ws.Cells.Item(86, 8).Value = 7239.1055
$ws.Cells.Item(86, 9).Value = 2085.7778
$ws.Cells.Item(86, 11).Value = 2085.7778
$ws.Cells.Item(86, 13).Value = -962.7777999999998

$ws.Cells.Item(89, 8).Value = 7239.1055
$ws.Cells.Item(89, 9).Value = 2085.7778
$ws.Cells.Item(89, 11).Value = 10428.889
$ws.Cells.Item(89, 13).Value = -4812.888999999999

$ws.Cells.Item(94, 8).Value = 5000
$ws.Cells.Item(94, 9).Value = 5000
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 5000
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 13).Value = -4549
$ws.Cells.Item(94, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 75000
$ws.Cells.Item(132, 10).Value = 75000
$ws.Cells.Item(132, 12).Value = 75000
$ws.Cells.Item(132, 14).Value = -85120

$ws.Cells.Item(134, 8).Value = 1684.6046
$ws.Cells.Item(134, 9).Value = 1708.0476
$ws.Cells.Item(134, 10).Value = 700
$ws.Cells.Item(134, 11).Value = 5124.142800000001
$ws.Cells.Item(134, 12).Value = 2100
$ws.Cells.Item(134, 13).Value = -2589.142800000001
$ws.Cells.Item(134, 14).Value = -7170

$ws.Cells.Item(135, 8).Value = 76476.8
$ws.Cells.Item(135, 10).Value = 76476.8
$ws.Cells.Item(135, 12).Value = 76476.8
$ws.Cells.Item(135, 14).Value = -86616.8

$ws.Cells.Item(137, 8).Value = 59500
$ws.Cells.Item(137, 10).Value = 59500
$ws.Cells.Item(137, 12).Value = 59500
$ws.Cells.Item(137, 14).Value = -69700

$ws.Cells.Item(138, 8).Value = 54999.332
$ws.Cells.Item(138, 10).Value = 54999.332
$ws.Cells.Item(138, 12).Value = 54999.332
$ws.Cells.Item(138, 14).Value = -65279.332

$ws.Cells.Item(139, 8).Value = 99985
$ws.Cells.Item(139, 10).Value = 99985
$ws.Cells.Item(139, 12).Value = 99985
$ws.Cells.Item(139, 14).Value = -110265

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 1094.2307
$ws.Cells.Item(122, 9).Value = 1093.75
$ws.Cells.Item(122, 10).Value = 1100
$ws.Cells.Item(122, 11).Value = 3281.25
$ws.Cells.Item(122, 12).Value = 3300
$ws.Cells.Item(122, 13).Value = -831.25
$ws.Cells.Item(122, 14).Value = -8200

$ws.Cells.Item(132, 8).Value = 3597.926
$ws.Cells.Item(132, 9).Value = 3514.5
$ws.Cells.Item(132, 10).Value = 3965
$ws.Cells.Item(132, 11).Value = 10543.5
$ws.Cells.Item(132, 12).Value = 11895
$ws.Cells.Item(132, 13).Value = -8013.5
$ws.Cells.Item(132, 14).Value = -16955

$ws.Cells.Item(134, 8).Value = 14634.241
$ws.Cells.Item(134, 9).Value = 5415.72
$ws.Cells.Item(134, 11).Value = 16247.16
$ws.Cells.Item(134, 13).Value = -13712.16

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 58
$ws.Cells.Item(17, 9).Value = 58
$ws.Cells.Item(17, 11).Value = 174
$ws.Cells.Item(17, 13).Value = -5

$ws.Cells.Item(34, 8).Value = 4631820.5
$ws.Cells.Item(34, 9).Value = 48.25
$ws.Cells.Item(34, 10).Value = 6947706.5
$ws.Cells.Item(34, 11).Value = 144.75
$ws.Cells.Item(34, 12).Value = 20843119.5
$ws.Cells.Item(34, 13).Value = -60.75
$ws.Cells.Item(34, 14).Value = -20843287.5

$ws.Cells.Item(39, 8).Value = 1936.875
$ws.Cells.Item(39, 10).Value = 1936.875
$ws.Cells.Item(39, 12).Value = 5810.625
$ws.Cells.Item(39, 14).Value = -6398.625

$ws.Cells.Item(55, 8).Value = 2261
$ws.Cells.Item(55, 10).Value = 3668.3333
$ws.Cells.Item(55, 12).Value = 11004.9999
$ws.Cells.Item(55, 14).Value = -11358.9999

$ws.Cells.Item(94, 8).Value = 5158.4165
$ws.Cells.Item(94, 10).Value = 5943.6665
$ws.Cells.Item(94, 12).Value = 17830.9995
$ws.Cells.Item(94, 14).Value = -19182.9995

$ws.Cells.Item(134, 8).Value = 11338.571
$ws.Cells.Item(134, 9).Value = 10067.272
$ws.Cells.Item(134, 10).Value = 16000
$ws.Cells.Item(134, 11).Value = 30201.816
$ws.Cells.Item(134, 12).Value = 48000
$ws.Cells.Item(134, 13).Value = -25131.816
$ws.Cells.Item(134, 14).Value = -58140

$ws.Cells.Item(137, 8).Value = 3621.5
$ws.Cells.Item(137, 10).Value = 5058.4443
$ws.Cells.Item(137, 12).Value = 15175.3329
$ws.Cells.Item(137, 14).Value = -25375.3329

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(41, 8).Value = 2249.75
$ws.Cells.Item(41, 9).Value = 2249.75
$ws.Cells.Item(41, 10).Value = 0
$ws.Cells.Item(41, 11).Value = 2249.75
$ws.Cells.Item(41, 12).Value = 0
$ws.Cells.Item(41, 13).Value = -1894.75
$ws.Cells.Item(41, 14).ClearContents()

$ws.Cells.Item(102, 8).Value = 37038532
$ws.Cells.Item(102, 9).Value = 1116.9445
$ws.Cells.Item(102, 11).Value = 1116.9445
$ws.Cells.Item(102, 13).Value = 505.0554999999999

$ws.Cells.Item(132, 8).Value = 4508.0835
$ws.Cells.Item(132, 9).Value = 4177.4443
$ws.Cells.Item(132, 10).Value = 5500
$ws.Cells.Item(132, 11).Value = 12532.3329
$ws.Cells.Item(132, 12).Value = 16500
$ws.Cells.Item(132, 13).Value = -10002.3329
$ws.Cells.Item(132, 14).Value = -21560

$ws.Cells.Item(141, 8).Value = 61905.6
$ws.Cells.Item(141, 10).Value = 61905.6
$ws.Cells.Item(141, 12).Value = 61905.6
$ws.Cells.Item(141, 14).Value = -72265.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(20, 8).Value = 17854.715
$ws.Cells.Item(20, 10).Value = 19596.6
$ws.Cells.Item(20, 12).Value = 19596.6
$ws.Cells.Item(20, 14).Value = -20048.6

$ws.Cells.Item(40, 8).Value = 7619
$ws.Cells.Item(40, 9).Value = 5316.3335
$ws.Cells.Item(40, 11).Value = 5316.3335
$ws.Cells.Item(40, 13).Value = -5180.3335

$ws.Cells.Item(136, 8).Value = 4345.3716
$ws.Cells.Item(136, 9).Value = 3522
$ws.Cells.Item(136, 11).Value = 10566
$ws.Cells.Item(136, 13).Value = -8016

$ws.Cells.Item(140, 8).Value = 92990
$ws.Cells.Item(140, 10).Value = 92990
$ws.Cells.Item(140, 12).Value = 92990
$ws.Cells.Item(140, 14).Value = -103350

$ws.Cells.Item(141, 8).Value = 164992.5
$ws.Cells.Item(141, 10).Value = 164992.5
$ws.Cells.Item(141, 12).Value = 164992.5
$ws.Cells.Item(141, 14).Value = -175352.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(140, 8).Value = 68755.8
$ws.Cells.Item(140, 10).Value = 68755.8
$ws.Cells.Item(140, 12).Value = 68755.8
$ws.Cells.Item(140, 14).Value = -79115.8
